$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.483.57'
$ws.Range("E2").Value = '  -1.65%  '
$ws.Range("D3").Value = '1.850.03'
$ws.Range("E3").Value = '  -1.57%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = '  +0.25%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '260.63'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  -8.22%  '
$ws.Range("E6").Value = '  +0.15%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5171'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = '  -1.50%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3239'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = '  -8.72%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06773'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  -4.44%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.98'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  -7.51%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7711'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  -6.83%  '
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07735'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = '  -0.36%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.865.36'
$ws.Range("E13").Value = '  -0.55%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '88.69'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = '  -1.80%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.038'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = '  -3.57%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.001'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = '  +0.31%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.13'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = '  -2.90%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9999'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = '  +0.13%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007907'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  -3.47%  '
$ws.Range("D20").Value = '26.544.11'
$ws.Range("E20").Value = '  -1.50%  '
$ws.Range("D21").Value = '2.129.62'
$ws.Range("E21").Value = '  +0.80%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.543'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = '  -5.33%  '
$ws.Range("E23").Value = '  -6.73%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.935'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  -5.25%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.345'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  -3.57%  '
$ws.Range("E26").Value = '  -1.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.658'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  -0.54%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.01'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = '  -2.98%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '111.59'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = '  -0.71%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.219'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  -4.89%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.177'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = '  -5.08%  '
$ws.Range("E32").Value = '  -1.52%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04811'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  -2.81%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.134'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  -4.69%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.850'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = '  -0.82%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6919'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  -8.54%  '
$ws.Range("E37").Value = '  -5.72%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01794'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  -5.20%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.214'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  -9.17%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.4909'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  -8.41%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '113.08'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = '  -3.50%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9039'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  -8.56%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.171'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  -2.66%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9998'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  +0.16%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.797'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  -5.91%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4219'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  -9.73%  '
$ws.Range("E47").Value = '  -8.11%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.168'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = '  -3.74%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05893'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = '  -0.83%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '35.43'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  -4.11%  '
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.421'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  -7.03%  '
